$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "abibangbranfffdon855@gmail.com"
$ws.Range("A5").Value = "rrrr@gmail.com"
$ws.Range("A6").Value = "a@gmail.com"
$ws.Range("A7").Value = "bbb@gmail.com"
$ws.Range("A8").Value = "ggffg@gmail.com"
